$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.4786585
$ws.Range("H2").Value = 38.957317
$ws.Range("I2").Value = 0.01644248566400343
$ws.Range("J2").Value = 0.01108359890151296
$ws.Range("Q2").Value = 5.427676247935667
$ws.Range("R2").Value = 32.566057487614
$ws.Range("S2").Value = 0.01644248566400343
$ws.Range("T2").Value = 0.01108359890151296

# Row 3
$ws.Range("I3").Value = 0.3309029145291901
$ws.Range("J3").Value = 0.3345840089140918
$ws.Range("S3").Value = 0.3309029145291901
$ws.Range("T3").Value = 0.3345840089140918

# Row 4
$ws.Range("G4").Value = 375.1018676666667
$ws.Range("H4").Value = 1125.305603
$ws.Range("I4").Value = 0.3166340783504202
$ws.Range("J4").Value = 0.3201564405802684
$ws.Range("Q4").Value = 104.5211351536696
$ws.Range("R4").Value = 940.6902163830262
$ws.Range("S4").Value = 0.3166340783504202
$ws.Range("T4").Value = 0.3201564405802684

# Row 5
$ws.Range("G5").Value = 19.6220475
$ws.Range("H5").Value = 39.244095
$ws.Range("I5").Value = 0.01656352436781744
$ws.Range("J5").Value = 0.01116518902553968
$ws.Range("Q5").Value = 5.467631210415001
$ws.Range("R5").Value = 32.80578726249001
$ws.Range("S5").Value = 0.01656352436781744
$ws.Range("T5").Value = 0.01116518902553968

# Row 6
$ws.Range("G6").Value = 133.4172743333333
$ws.Range("H6").Value = 400.2518229999999
$ws.Range("I6").Value = 0.1126212886044614
$ws.Range("J6").Value = 0.1138741321875775
$ws.Range("Q6").Value = 37.17636771358512
$ws.Range("R6").Value = 334.587309422266
$ws.Range("S6").Value = 0.1126212886044614
$ws.Range("T6").Value = 0.1138741321875775

# Row 7
$ws.Range("G7").Value = 245.0287756666667
$ws.Range("H7").Value = 735.086327
$ws.Range("I7").Value = 0.2068357084841073
$ws.Range("J7").Value = 0.2091366303910096
$ws.Range("Q7").Value = 68.27661492944823
$ws.Range("R7").Value = 614.489534365034
$ws.Range("S7").Value = 0.2068357084841073
$ws.Range("T7").Value = 0.2091366303910096

$wb.Save()
